$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "Total:"
$ws.Range("B19").Formula = "=SUM(Sheet2!D2,Sheet2!D11)"
